$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New body text for the "e051 Crew Actions - Tank Movement" entry (shared string added to sharedStrings.xml)
$e051Body = @'
<Bold>e051 Crew Actions - Tank Movement</Bold> 
<InlineUIContainer><Button Content='r4.74.1' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>   
<InlineUIContainer><Button Content='r11.0' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>  
<LineBreak/><LineBreak/>
Resolve movement per the 
<InlineUIContainer><Button Content='Movement' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> Tables.
<LineBreak/><LineBreak/>
<Underline>Modifiers:</Underline><LineBreak/>
'@

# Insert a new row above the current row 61 ("e100" row), shifting rows 61-70 down to 62-71,
# then populate the new row with the e051 entry.
$ws.Rows("61:61").Insert()

$ws.Range("A61").Value = "e051"
$ws.Range("B61").Value = $e051Body
$ws.Rows("61:61").RowHeight = 120

# Update the view to reflect scrolling to / selecting the newly added row, like a user would
# after typing the new entry.
$win = $excel.ActiveWindow
$win.ScrollRow = 56
$win.ScrollColumn = 1
$ws.Range("B61").Select()
